$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new expense rows (row 6 and row 7) matching the existing data pattern
$ws.Range("A6").Value = "Despesa"
$ws.Range("B6").Value = "SERVIÇOS"
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = "26/01/2025"

$ws.Range("A7").Value = "Despesa"
$ws.Range("B7").Value = "SERVIÇOS"
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = "26/01/2025"
